# Add a second Login test case (LOG-TC-02) to the "Login Test" sheet,
# mirroring the structure/format of the existing LOG-TC-01 row (row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy row 4's formatting down into row 5 so the new row matches the
#     existing look (borders, fills, alignment, wrap, fonts, etc.) ---
$ws.Range("B4:L4").Copy()
$ws.Range("B5:L5").PasteSpecial(-4122)

# --- Fill in the new test case's data ---
$ws.Range("B5").Value = "LOG-02"
$ws.Range("C5").Value = "Verify system displays error when password is incorrect"
$ws.Range("D5").Value = "LOG-TC-02"
$ws.Range("E5").Value = "Login Module"
$ws.Range("F5").Value = "LOG-TS-02"
$ws.Range("G5").Value = "Login with incorrect password"
$ws.Range("H5").Value = "Email already registered"
$ws.Range("I5").Value = "1. Open login page `n2. Enter registered email `n3. Enter wrong password `n4. Click Login"
$ws.Range("J5").Value = "user1@gmail.com`nWrongPass"
$ws.Range("K5").Value = "Error message displayed"
$ws.Range("L5").Value = "High"

# --- Add the mailto hyperlink on J5 (test data), matching the pattern on J4 ---
$ws.Hyperlinks.Add($ws.Range("J5"), "mailto:user1@gmail.com%0aWrongPass")

# Re-apply J4's formatting onto J5 (Hyperlinks.Add resets to the default
# "Hyperlink" style) so it keeps the bordered/centered/wrapped look.
$ws.Range("J4").Copy()
$ws.Range("J5").PasteSpecial(-4122)

# --- Row height for the new row (wraps to a few lines like row 4) ---
$ws.Rows.Item(5).RowHeight = 78.75

# --- Widen column C slightly to fit the new longer description ---
$ws.Columns.Item(3).ColumnWidth = 49.67

# --- Update the sheet view: drop the frozen/scrolled "E1" top-left cell
#     and move the active selection to A6 ---
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A6").Select()
